$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.089808712187860351
$ws.Range("A2").Value = -0.0099999997733242196
$ws.Range("A3").Value = -0.0089999997750069838
$ws.Range("A4").Value = 0.283996363610445
$ws.Range("A5").Value = -0.0059999997824515816
$ws.Range("A6").Value = -0.017893757874766436
$ws.Range("A7").Value = -0.019999999733270712
$ws.Range("A8").Value = -0.019999999732319473
$ws.Range("A9").Value = -0.0059999997730608712
$ws.Range("A10").Value = -0.031113678696776503
$ws.Range("A11").Value = -0.0044999997776358214
$ws.Range("A12").Value = 0.055676827615191549
$ws.Range("A13").Value = -0.0059999997722508525
$ws.Range("A14").Value = -0.011999999754212176
$ws.Range("A15").Value = 0.018950289939039244
$ws.Range("A16").Value = -0.0059999997709927477
$ws.Range("A17").Value = -0.0059999997699478058
$ws.Range("A18").Value = -0.0089999997608440907
$ws.Range("A19").Value = -0.0089999997760776829
$ws.Range("A20").Value = -0.0089999997739909077
$ws.Range("A21").Value = -0.0089999997736605053
$ws.Range("A22").Value = -0.0089999997733896109
$ws.Range("A23").Value = -0.062820653472729582
$ws.Range("A24").Value = -0.04199999967058865
$ws.Range("A25").Value = -0.041999999668608012
$ws.Range("A26").Value = -0.0059999997739588196
$ws.Range("A27").Value = -0.0059999997727482324
$ws.Range("A28").Value = -0.0059999997676802863
$ws.Range("A29").Value = -0.011999999746363343
$ws.Range("A30").Value = -0.019999999721099559
$ws.Range("A31").Value = -0.014999999733285918
$ws.Range("A32").Value = -0.020999999715142437
$ws.Range("A33").Value = -0.0059999997591599907
